# Apply cryptos list update (prices / 1h volume deltas) per diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.042.41"
$ws.Range("E2").Value = "  +2.49%  "
$ws.Range("D3").Value = "3.030.52"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.98"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.52"
$ws.Range("E6").Value = "  +6.56%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.025.96"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("E10").Value = "  +17.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("E11").Value = "  +1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.65"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").Value = "3.529.40"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").Value = "63.008.16"
$ws.Range("E17").Value = "  +2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.07"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("D19").Value = "3.029.19"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.79"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.695"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.53"
$ws.Range("E24").Value = "  +8.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.21"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("E26").Value = "  +7.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.37"
$ws.Range("E27").Value = "  +3.09%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("E30").Value = "  +10.92%  "
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.67"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +5.70%  "
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.89"
$ws.Range("E37").Value = "  +2.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("E38").Value = "  +10.28%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.131"
$ws.Range("E39").Value = "  +9.02%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.11"
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.53"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.07"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.19"
$ws.Range("E44").Value = "  +15.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "394.00"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0360"
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("D47").Value = "2.719.18"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.60"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.68"
$ws.Range("E49").Value = "  +14.00%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("E51").Value = "  +6.14%  "
